$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.766.28"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.717.75"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.79"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.46"
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.718.59"
$ws.Range("E7").Value = "  -1.55%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  +4.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.24"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.14"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000243"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.348.04"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.720.97"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.934.15"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.31"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.37"
$ws.Range("E20").Value = "  +8.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "488.43"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.25"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.727"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.86"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000142"
$ws.Range("E25").Value = "  +4.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  -1.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.33"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.37"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.77"
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.47"
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.864.94"
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.667.84"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.85"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.324"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "431.39"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.91"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.96"
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.46"
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.72"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.27"
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.758.43"
$ws.Range("E51").Value = "  -3.09%  "
